$d = $word.ActiveDocument

# 1. Title change - occurs twice (Heading1 at top, and bold run near the end).
#    Replace all occurrences in one pass.
$d.Content.Find.Execute(
    "Play Lucky Tanks for Free - Slot Game Review", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Lucky Tanks Free - Unique Slot Game with Total Player Control", 2)

# 2. "What we like" bullet list.
#    Insert a brand-new bullet ("Unique combination of slot machines and lotteries")
#    right before the existing "Total player control over gameplay" bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Total player control over gameplay`r") {
        $p.Range.InsertParagraphBefore()
        break
    }
}

# InsertParagraphBefore() leaves a fresh empty "List Bullet" paragraph just
# ahead of the "Total player control..." one - find it and fill in its text.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "List Bullet" -and $p.Range.Text -eq "`r") {
        $p.Range.Text = "Unique combination of slot machines and lotteries"
        break
    }
}

# Replace remaining bullets in "What we like".
$d.Content.Find.Execute(
    "Combination of slot machines and lotteries", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Military theme with explosions and military equipment", 2)

$d.Content.Find.Execute(
    "Valid chances of winning with fixed RTP", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Simple and intuitive design and graphics", 2)

# Remove the now-redundant "Well-developed graphics and military theme" bullet entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Well-developed graphics and military theme`r") {
        $p.Range.Delete()
        break
    }
}

# 3. "What we don't like" bullet list.
$d.Content.Find.Execute(
    "Simple gameplay may require patience", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Requires patience for gameplay", 2)

$d.Content.Find.Execute(
    "Limited number of usable weapons", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "No extra functions often expected in traditional slots", 2)

# 4. Italic summary paragraph near the end.
$d.Content.Find.Execute(
    "Read our review of Lucky Tanks, the slot game that combines slot machines and lotteries, and play it for free. Full game features and mechanics explained.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Lucky Tanks for free and experience a unique combination of slot machines and lotteries with total player control over gameplay.", 2)
